$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 145, shifting existing rows 145:266 down to 146:267.
$ws.Rows("145:145").Insert()

# Fill the newly inserted row 145 with the new data record.
$ws.Cells.Item(145, 1).Value = 5
$ws.Cells.Item(145, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(145, 3).Value = "Maule"
$ws.Cells.Item(145, 4).Value = 44741
$ws.Cells.Item(145, 5).Value = 7
$ws.Cells.Item(145, 6).Value = 100112009
$ws.Cells.Item(145, 7).Value = "Acelga"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 500
$ws.Cells.Item(145, 11).Value = 3000
$ws.Cells.Item(145, 12).Value = 3000
$ws.Cells.Item(145, 13).Value = 3000
$ws.Cells.Item(145, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(145, 15).Value = "Región del Maule"
$ws.Cells.Item(145, 16).Value = 750
$ws.Cells.Item(145, 17).Value = 4
$ws.Cells.Item(145, 18).Value = "Hortaliza"
